$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number, Report Covering the Week dates) ---
$ws.Range("A8").Value = "Volume 30   Number  18"
$ws.Range("C9").Value = "Report Covering the Week  5/1/2023  Through  5/7/2023"

# --- Row 22: Fel. Assault cells D22/E22 convert from "no data" placeholders to real numbers ---
$ws.Range("D22").NumberFormat = "#,##0"
$ws.Range("D22").Value = 3
$ws.Range("E22").NumberFormat = "#,##0.0;""-""#,##0.0"
$ws.Range("E22").Value = -100

# --- Crime Complaints table value updates (rows 14-29) ---
$ws.Range("L14").Value = -33.333333333333
$ws.Range("N14").Value = -81.818181818181
$ws.Range("F15").Value = 7
$ws.Range("H15").Value = 250
$ws.Range("I15").Value = 15
$ws.Range("K15").Value = 25
$ws.Range("L15").Value = 7.142857142857
$ws.Range("M15").Value = 150
$ws.Range("N15").Value = 36.363636363636
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 29
$ws.Range("G16").Value = 15
$ws.Range("H16").Value = 93.333333333333
$ws.Range("I16").Value = 103
$ws.Range("J16").Value = 87
$ws.Range("K16").Value = 18.390804597701
$ws.Range("L16").Value = 71.666666666666
$ws.Range("M16").Value = 19.767441860465
$ws.Range("N16").Value = -63.475177304964
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 20
$ws.Range("G17").Value = 30
$ws.Range("H17").Value = -10
$ws.Range("I17").Value = 127
$ws.Range("J17").Value = 134
$ws.Range("K17").Value = -5.223880597014
$ws.Range("L17").Value = 53.012048192771
$ws.Range("M17").Value = -6.617647058823
$ws.Range("N17").Value = -7.971014492753
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 2
$ws.Range("E18").Value = 0
$ws.Range("F18").Value = 15
$ws.Range("G18").Value = 15
$ws.Range("I18").Value = 66
$ws.Range("J18").Value = 52
$ws.Range("K18").Value = 26.923076923076
$ws.Range("L18").Value = 73.684210526315
$ws.Range("M18").Value = -19.512195121951
$ws.Range("N18").Value = -77.926421404682
$ws.Range("C19").Value = 4
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = -33.333333333333
$ws.Range("F19").Value = 31
$ws.Range("G19").Value = 23
$ws.Range("H19").Value = 34.782608695652
$ws.Range("I19").Value = 138
$ws.Range("J19").Value = 138
$ws.Range("K19").Value = 0
$ws.Range("L19").Value = 42.268041237113
$ws.Range("M19").Value = 56.818181818181
$ws.Range("N19").Value = 8.661417322834
$ws.Range("C20").Value = 4
$ws.Range("D20").Value = 4
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 27
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = 80
$ws.Range("I20").Value = 121
$ws.Range("J20").Value = 83
$ws.Range("K20").Value = 45.783132530120
$ws.Range("L20").Value = 152.083333333333
$ws.Range("M20").Value = 348.148148148148
$ws.Range("N20").Value = -9.701492537313
$ws.Range("C21").Value = 25
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = 4.166666666666
$ws.Range("F21").Value = 136
$ws.Range("G21").Value = 101
$ws.Range("H21").Value = 34.653465346534
$ws.Range("I21").Value = 572
$ws.Range("J21").Value = 507
$ws.Range("K21").Value = 12.820512820512
$ws.Range("L21").Value = 66.763848396501
$ws.Range("M21").Value = 34.588235294117
$ws.Range("N21").Value = -42.914171656686
$ws.Range("F22").Value = 1
$ws.Range("G22").Value = 3
$ws.Range("H22").Value = -66.666666666666
$ws.Range("J22").Value = 10
$ws.Range("K22").Value = -40
$ws.Range("M22").Value = -14.285714285714
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 21
$ws.Range("E24").Value = -14.285714285714
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 58
$ws.Range("H24").Value = 17.241379310344
$ws.Range("I24").Value = 266
$ws.Range("J24").Value = 292
$ws.Range("K24").Value = -8.904109589041
$ws.Range("L24").Value = 59.281437125748
$ws.Range("M24").Value = 49.438202247191
$ws.Range("C25").Value = 9
$ws.Range("D25").Value = 14
$ws.Range("E25").Value = -35.714285714285
$ws.Range("F25").Value = 44
$ws.Range("G25").Value = 45
$ws.Range("H25").Value = -2.222222222222
$ws.Range("I25").Value = 165
$ws.Range("J25").Value = 174
$ws.Range("K25").Value = -5.172413793103
$ws.Range("L25").Value = 25
$ws.Range("M25").Value = -14.0625
$ws.Range("C26").Value = 4
$ws.Range("F26").Value = 11
$ws.Range("G26").Value = 4
$ws.Range("H26").Value = 175
$ws.Range("I26").Value = 21
$ws.Range("K26").Value = 0
$ws.Range("L26").Value = -25
$ws.Range("C27").Value = 2
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 0
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 120
$ws.Range("I27").Value = 36
$ws.Range("J27").Value = 43
$ws.Range("K27").Value = -16.279069767441
$ws.Range("L27").Value = -5.263157894736
$ws.Range("L28").Value = 50
$ws.Range("L29").Value = 25
